$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the coordinates for object 50433 (row 9: A=50433, B=z, C=RA, D=Dec)
$ws.Range("C9").Value = 157.75800000000001
$ws.Range("D9").Value = 30.861999999999998

# Update the active cell selection to reflect where the user left off editing
$ws.Range("D12").Select()
